{"js": "// The template has a paragraph containing a Word field whose field code is\n// \" m:x \" (fldChar begin / instrText \" \" \"m\" \":x\" \" \" / fldChar end). This\n// edit turns that field into plain literal text runs spelling out\n// \"{\", \"m\", \":x\", \"}\" so the paragraph reads \"{m:x}\" as plain text instead\n// of a field (matching the move to TokenIteratorFieldRewriterSplit parsing\n// of literal \"{m:x}\" template markers instead of Word fields).\n\nconst replacementOoxml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' +\n  '<w:p>' +\n  '<w:r><w:t>{</w:t></w:r>' +\n  '<w:r><w:t>m</w:t></w:r>' +\n  '<w:r><w:t>:x</w:t></w:r>' +\n  '<w:r><w:t xml:space=\"preserve\">}</w:t></w:r>' +\n  '</w:p>' +\n  '</w:body></w:document>' +\n  '</pkg:xmlData></pkg:part></pkg:package>';\n\nlet targetRange = null;\n\n// Preferred approach: use the document Fields collection to find the field\n// whose code is \"m:x\" (ignoring surrounding whitespace) and target the\n// paragraph that hosts it.\nconst fields = context.document.body.fields;\nfields.load(\"items\");\nawait context.sync();\n\nif (fields.items.length > 0) {\n  for (let i = 0; i < fields.items.length; i++) {\n    fields.items[i].load(\"code\");\n  }\n  await context.sync();\n\n  let targetField = null;\n  for (let i = 0; i < fields.items.length; i++) {\n    const code = (fields.items[i].code || \"\").replace(/\\s+/g, \"\");\n    if (code === \"m:x\") {\n      targetField = fields.items[i];\n      break;\n    }\n  }\n\n  if (targetField) {\n    const fieldParagraphs = targetField.result.paragraphs;\n    fieldParagraphs.load(\"items\");\n    await context.sync();\n    if (fieldParagraphs.items.length > 0) {\n      targetRange = fieldParagraphs.items[0].getRange(\"Whole\");\n    }\n  }\n}\n\n// Fallback: scan every paragraph's OOXML looking for the instrText field\n// code markup containing \"m:x\" (used if the Fields API above is\n// unavailable or did not locate the field for any reason).\nif (!targetRange) {\n  const paragraphs = context.document.body.paragraphs;\n  paragraphs.load(\"items\");\n  await context.sync();\n\n  const ooxmlResults = paragraphs.items.map((p) => p.getRange(\"Whole\").getOoxml());\n  await context.sync();\n\n  for (let i = 0; i < paragraphs.items.length; i++) {\n    const xml = ooxmlResults[i].value;\n    if (xml.indexOf(\"w:instrText\") !== -1 && /:\\s*x/.test(xml) && />\\s*m\\s*</.test(xml)) {\n      targetRange = paragraphs.items[i].getRange(\"Whole\");\n      break;\n    }\n  }\n}\n\nif (!targetRange) {\n  throw new Error(\"Could not locate the paragraph containing the 'm:x' field.\");\n}\n\ntargetRange.insertOoxml(replacementOoxml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# The template has a paragraph containing a Word field whose field code is\n# \" m:x \" (fldChar begin / instrText \" \", \"m\", \":x\", \" \" / fldChar end).\n# This edit turns that field into plain literal text runs spelling out\n# \"{\", \"m\", \":x\", \"}\" so the paragraph reads \"{m:x}\" as plain text instead\n# of a field.\n\n$d = $word.ActiveDocument\n\n# Locate the field whose code (trimmed of whitespace) is \"m:x\".\n$targetField = $null\nforeach ($f in $d.Fields) {\n    $code = $f.Code.Text\n    $trimmed = ($code -replace '\\s+', '')\n    if ($trimmed -eq \"m:x\") {\n        $targetField = $f\n        break\n    }\n}\n\nif ($targetField -eq $null) {\n    throw \"Could not find the field with code 'm:x'.\"\n}\n\n$codeStart = $targetField.Code.Start\n$codeEnd = $targetField.Code.End\n\n# Find the paragraph that contains the field code range.\n$targetParagraph = $null\nforeach ($p in $d.Paragraphs) {\n    $pr = $p.Range\n    if ($codeStart -ge $pr.Start -and $codeEnd -le $pr.End) {\n        $targetParagraph = $p\n        break\n    }\n}\n\nif ($targetParagraph -eq $null) {\n    throw \"Could not find the paragraph hosting the 'm:x' field.\"\n}\n\n$range = $targetParagraph.Range\n\n$xml = \"<pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'>\" +\n       \"<pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'>\" +\n       \"<pkg:xmlData>\" +\n       \"<w:document xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>\" +\n       \"<w:body>\" +\n       \"<w:p>\" +\n       \"<w:r><w:t>{</w:t></w:r>\" +\n       \"<w:r><w:t>m</w:t></w:r>\" +\n       \"<w:r><w:t>:x</w:t></w:r>\" +\n       \"<w:r><w:t xml:space='preserve'>}</w:t></w:r>\" +\n       \"</w:p>\" +\n       \"</w:body></w:document>\" +\n       \"</pkg:xmlData></pkg:part></pkg:package>\"\n\n$range.InsertXML($xml)\n"}
